$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Brasil", "01/01/2015", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.72503277584889, $null, $true),
    @("Brasil", "01/01/2016", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.97489130555118, $null, $true),
    @("Brasil", "01/01/2017", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 13.68232552634555, $null, $false),
    @("Brasil", "01/01/2018", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 13.11379688250342, $null, $false),
    @("Brasil", "01/01/2019", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 11.94262203013007, $null, $false),
    @("Brasil", "01/01/2020", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 11.79791843848013, $null, $false),
    @("Brasil", "01/01/2021", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.19148357379749, $null, $false),
    @("Brasil", "01/01/2022", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.75457740723266, $null, $false),
    @("Brasil", "01/01/2023", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.14424902108152, $null, $false),
    @("Brasil", "01/01/2024", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 13.93153959937777, $null, $false),
    @("Brasil", "01/01/2025", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 7.183472099432135, $null, $false),
    @("Nordeste", "01/01/2015", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 16.71215007526446, $null, $true),
    @("Nordeste", "01/01/2016", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 16.45884553085288, $null, $true),
    @("Nordeste", "01/01/2017", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 13.89967984937273, $null, $false),
    @("Nordeste", "01/01/2018", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 14.76660569953508, $null, $false),
    @("Nordeste", "01/01/2019", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.63533242853311, $null, $false),
    @("Nordeste", "01/01/2020", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.47263083599161, $null, $false),
    @("Nordeste", "01/01/2021", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.44979916779812, $null, $false),
    @("Nordeste", "01/01/2022", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.58606062940625, $null, $false),
    @("Nordeste", "01/01/2023", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.20043826255992, $null, $false),
    @("Nordeste", "01/01/2024", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 13.78132210227987, $null, $false),
    @("Nordeste", "01/01/2025", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 7.278363804704043, $null, $false),
    @("Sergipe", "01/01/2015", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 15.33841275398043, 9, $true),
    @("Sergipe", "01/01/2016", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 12.47132489276449, 12, $true),
    @("Sergipe", "01/01/2017", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 10.1007147584733, 20, $false),
    @("Sergipe", "01/01/2018", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 11.499762104158, 15, $false),
    @("Sergipe", "01/01/2019", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 8.308314143012621, 22, $false),
    @("Sergipe", "01/01/2020", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 6.468216048765174, 24, $false),
    @("Sergipe", "01/01/2021", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 6.200057468118876, 24, $false),
    @("Sergipe", "01/01/2022", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 5.089403861161063, 26, $false),
    @("Sergipe", "01/01/2023", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 6.648580843587086, 23, $false),
    @("Sergipe", "01/01/2024", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 8.143716976183596, 21, $false),
    @("Sergipe", "01/01/2025", "Morte no trânsito ou em decorrência dele (exceto homicídio doloso)", 3.648542531640742, 25, $false),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    if ($null -eq $vals[4]) {
        $ws.Cells.Item($row, 5).Value = ""
    } else {
        $ws.Cells.Item($row, 5).Value = $vals[4]
    }
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
